# Refresh the "cryptos" price table (rows 2-51) with the latest scrape.
# Price (col D) and Volume(1h) (col E) are text cells (source data uses
# locale-style separators such as "30.577.53", so Excel cannot treat them
# as real numbers) and a couple of rows got reordered (their Coin/Link
# swapped) because the ranking shifted between runs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as TEXT (avoids Excel auto-numeric coercion
# for numeric-looking strings like "245.38" or "0.9998"), while leaving
# the cell style/number-format exactly as it was before (General, no custom xf).
function Set-TextValue($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" '30.576.90'
Set-TextValue $ws "E2" '  -0.09%  '
Set-TextValue $ws "D3" '1.921.95'
Set-TextValue $ws "E3" '  -0.04%  '
Set-TextValue $ws "E4" '  -0.02%  '
Set-TextValue $ws "D5" '245.38'
Set-TextValue $ws "E5" '  -1.11%  '
Set-TextValue $ws "E6" '  -0.05%  '
Set-TextValue $ws "D7" '0.4830'
Set-TextValue $ws "E7" '  +1.99%  '
Set-TextValue $ws "D8" '0.2897'
Set-TextValue $ws "E8" '  -0.26%  '
Set-TextValue $ws "D9" '0.06796'
Set-TextValue $ws "E9" '  -0.78%  '
Set-TextValue $ws "D10" '112.19'
Set-TextValue $ws "E10" '  +6.47%  '
Set-TextValue $ws "D11" '19.48'
Set-TextValue $ws "E11" '  +6.01%  '
Set-TextValue $ws "D12" '1.913.53'
Set-TextValue $ws "E12" '  -0.43%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws "D13" '0.07569'
Set-TextValue $ws "E13" '  -1.91%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws "D14" '5.463'
Set-TextValue $ws "E14" '  +2.10%  '
Set-TextValue $ws "D15" '0.6747'
Set-TextValue $ws "E15" '  +0.73%  '
Set-TextValue $ws "D16" '293.20'
Set-TextValue $ws "E16" '  +1.35%  '
Set-TextValue $ws "D17" '30.579.63'
Set-TextValue $ws "E17" '  -0.09%  '
Set-TextValue $ws "D18" '0.000007662'
Set-TextValue $ws "E18" '  +0.25%  '
Set-TextValue $ws "D19" '13.03'
Set-TextValue $ws "E19" '  +0.67%  '
Set-TextValue $ws "D20" '0.9998'
Set-TextValue $ws "D21" '5.509'
Set-TextValue $ws "E21" '  -0.60%  '
Set-TextValue $ws "D22" '2.164.76'
Set-TextValue $ws "E22" '  -0.54%  '
Set-TextValue $ws "D23" '1.001'
Set-TextValue $ws "E23" '  +0.02%  '
Set-TextValue $ws "D24" '6.453'
Set-TextValue $ws "E24" '  +0.00%  '
Set-TextValue $ws "D25" '9.485'
Set-TextValue $ws "E25" '  -0.29%  '
Set-TextValue $ws "D26" '167.03'
Set-TextValue $ws "E26" '  -0.42%  '
Set-TextValue $ws "D27" '20.32'
Set-TextValue $ws "E27" '  -2.06%  '
Set-TextValue $ws "D28" '2.101'
Set-TextValue $ws "E28" '  -1.05%  '
Set-TextValue $ws "D29" '0.1066'
Set-TextValue $ws "E29" '  -0.56%  '
Set-TextValue $ws "D30" '1.441'
Set-TextValue $ws "E30" '  +2.34%  '
Set-TextValue $ws "D31" '4.140'
Set-TextValue $ws "E31" '  -0.99%  '
Set-TextValue $ws "D32" '4.059'
Set-TextValue $ws "E32" '  +0.12%  '
Set-TextValue $ws "D33" '0.04944'
Set-TextValue $ws "E33" '  -1.45%  '
Set-TextValue $ws "D34" '0.7359'
Set-TextValue $ws "E34" '  +0.25%  '
Set-TextValue $ws "D35" '1.139'
Set-TextValue $ws "E35" '  -0.68%  '
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws "D36" '2.714'
Set-TextValue $ws "E36" '  -0.48%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws "D37" '0.02030'
Set-TextValue $ws "E37" '  -2.07%  '
Set-TextValue $ws "E38" '  +0.08%  '
Set-TextValue $ws "D39" '2.021'
Set-TextValue $ws "E39" '  -1.07%  '
Set-TextValue $ws "D40" '109.69'
Set-TextValue $ws "E40" '  -1.42%  '
Set-TextValue $ws "D41" '0.4438'
Set-TextValue $ws "E41" '  +0.22%  '
Set-TextValue $ws "D42" '0.8703'
Set-TextValue $ws "E42" '  -0.66%  '
Set-TextValue $ws "D43" '5.841'
Set-TextValue $ws "E43" '  -0.90%  '
Set-TextValue $ws "E44" '  +0.01%  '
Set-TextValue $ws "D45" '69.23'
Set-TextValue $ws "E45" '  +2.06%  '
Set-TextValue $ws "D46" '7.245'
Set-TextValue $ws "E46" '  -0.80%  '
Set-TextValue $ws "D47" '48.79'
Set-TextValue $ws "E47" '  +1.76%  '
Set-TextValue $ws "D48" '9.259'
Set-TextValue $ws "E48" '  -1.43%  '
Set-TextValue $ws "D49" '0.1232'
Set-TextValue $ws "E49" '  -0.77%  '
Set-TextValue $ws "D50" '34.86'
Set-TextValue $ws "E50" '  -0.35%  '
Set-TextValue $ws "D51" '0.2497'
Set-TextValue $ws "E51" '  -0.21%  '
